$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header in A1 and clear out the "Model" column (E)
$ws.Range("A1").Value = "Sample_Basename"
$ws.Range("E1:E3").ClearContents()

# Update selection to E1
$ws.Range("E1").Select()
